$d = $word.ActiveDocument

$newLines = @(
    "22/09/23",
    "Began Version 0.2.0 of the program",
    "2/10/23",
    "Completed Trialling and Testing Version 0.2.0 of the program",
    "6/10/23",
    "Completed Version 0.2.0 of the program",
    "Created .zip file containing executable for the program"
)

foreach ($line in $newLines) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $line
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
